# Applies the "Dynamic block now runs as loop instead of hard coding each
# metabolite/statevar" change to the "Initial & Flow Concentrations" sheet.
#
# Previously columns A (Initial Condition) and B (Flow Concentration) held the
# values for every metabolite/state-variable, hard-coded row by row. The new
# layout keeps A/B zeroed out (they are now populated dynamically, at runtime,
# by the simulation loop) and relocates the old hard-coded numbers to columns
# F (old initial condition) / G (old flow concentration) so they remain
# available as a reference/lookup table for the loop to read from.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Initial & Flow Concentrations")

# Row 1 is a header/parameter row (D = dilution rate); only the first value
# flips from 0 to 1.
$ws.Range("A1").Value2 = 1

# Rows 2-31 each describe one metabolite / state variable.
for ($r = 2; $r -le 31; $r++) {

    $oldB = $ws.Cells.Item($r, 2).Value2   # old "Flow Concentration"
    if ($null -eq $oldB) { $oldB = 0 }

    # Only rows 2-5 still had their old initial-condition value sitting in
    # column A (rows 6+ were already left blank there historically), so only
    # those rows need that value relocated into column F.
    if ($r -ge 2 -and $r -le 5) {
        $oldA = $ws.Cells.Item($r, 1).Value2   # old "Initial Condition"
        if ($null -eq $oldA) { $oldA = 0 }

        $ws.Cells.Item($r, 6).Value2 = $oldA   # F column
        $ws.Cells.Item($r, 1).Value2 = 0       # A becomes the dynamic placeholder
        $ws.Cells.Item($r, 1).Style = "Normal" # no longer hand-formatted
    }

    $ws.Cells.Item($r, 7).Value2 = $oldB       # G column holds the relocated flow concentration

    # B now just holds the dynamic (loop-driven) placeholder value 0.
    $ws.Cells.Item($r, 2).Value2 = 0
    $ws.Cells.Item($r, 2).Style = "Normal"     # no longer hand-formatted
}

# The previously cached selection (K20) is no longer meaningful once the
# sheet is rebuilt by the loop, so clear it.
$ws.Range("A1").Select()
